# Adds the 4 new component rows (entretoises, connecteur femelle jack,
# chargeur 12v, ressorts) to the "Liste de composant" sheet, matching the
# rows appended in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Component names were typed first (column A, top to bottom) ...
$ws.Range("A6").Value = "entretoises"
$ws.Range("A7").Value = "connecteur femelle jack DC 5,5X2,1"
$ws.Range("A8").Value = "chargeur 12v"
$ws.Range("A9").Value = "ressorts (X5)"

# ... then the prices (column B) ...
$ws.Range("B7").Value = 3.49
$ws.Range("B8").Value = 11.97
$ws.Range("B9").Value = 6.86

# ... then the links (column C), matching the order new strings were
# appended to the shared-strings table in the source file.
$ws.Range("C7").Value = 'https://www.amazon.fr/BIlinli-DC-005-Connecteur-5-5x2-1mm-Barrel/dp/B082B7VQRY/ref=sr_1_6?__mk_fr_FR=%C3%85M%C3%85%C5%BD%C3%95%C3%91&dib=eyJ2IjoiMSJ9.zaqWM5kZHUhQQR3Ar-_dPFqrPcrjUAWIURpgQoX8U4pKqx9iJOyQmh8juWfIply3MGozRWA4PK5UpUi6uHCja_YsyH3aIb02I321vk5zIY0zZrQo0s-_rcP59bzlGbUn3jXrk-PYGWoNa-F-5LTEYbHK9u-Qu3S-paWdsLBG8eU1d6oI1D5kuIkblKNC9bjMXY-NbJPzpQ5Z6PnpOCU6Y2fHgQs9qCvxiMC2G95dktSBfqcqL0Oeq9AvS1Cm-QWJb5yTROig-u1Veq1XgvrIeydLv7gQLXi-ckNoXkAq_5k2Goki83N-cCKUCYsvGkQXMdChRscXAuSHN7zVTiB4I1VhOMp0eVQ3LzpFCiTdLgKQPikedw7iBBE7KbatGXVx2D_yQLIKvAJyr8gAgesbjE6H1z0uCGF7OwR8jupltClHrsQ4o_qMF9-yYjMGd9fD.GwQgSGeDNs8PSCiMyg35wu3cF7VkPHT0IVfdAYKmNRU&dib_tag=se&keywords=dc+jack+pcb&qid=1740840369&sr=8-6'
$ws.Range("C8").Value = 'https://www.amazon.fr/Adaptateur-Alimentation-dalimentation-adaptateurs-convertisseur/dp/B09KNC7SDQ/ref=sr_1_5?dib=eyJ2IjoiMSJ9.kh4ME1fgLwcNRdIka5gVXrYsbm_HulQVUeVuKNMClm8kezgqD7oLYXZEnvZhdBtjKVmlWR_ygZ47e2nU3gspl7XBdliMYf5fC_xpPFQp4p4zdewgWqvd9Ztq_yhXsBIQPHK9tyvEENLNJBp6HMmjhLIoaLHOeYguf6Pe40BN_MimRcglxRfbslwaohYN5KlC7XeBYpvazD_MMvZSYMF1LUUodypRUWlFoDmgu4qjOvG3QU4wBi_TmTjXj5H_GzkbjUMp1ByKSg3463KvWQs4A-8MfTfp8Jt-j8qL7jwZO9wv2mbRI3syJ0x0MFgpKITtZXvGRxiolqiw0f7XFFc9P2wwoA0Xi03zowoz5NAV5y8r04eN34F3d36OsWB0umprTbaXlaY1RsjP_8t2GAV269poIuAMs7Qx47t51wI3Hhrq5XRGMpno8_4ibh5r6gwe.OO38cFMDgwd4r9Wqpo2xCEjt3winmd1ijKS_Ee3Ozqc&dib_tag=se&keywords=transformateur%2B12v&qid=1740840232&sr=8-5&th=1'
$ws.Range("C9").Value = 'https://fr.rs-online.com/web/p/ressorts-de-compression/0751512?searchId=1b93648f-6bca-46c3-8c8a-d5bf4eda4f0b&gb=s'

# Column widths grew slightly after the edit (re-measured by Excel once the
# new, longer component names/links were entered).
$ws.Columns.Item(1).ColumnWidth = 30.21875
$ws.Columns.Item(2).ColumnWidth = 10.109375
$ws.Columns.Item(3).ColumnWidth = 39.33203125

# Move the active selection down past the newly entered data, like a user
# would after typing the last row and pressing Enter a few times.
$ws.Range("C19").Select()
